$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 17.10369991047524
$ws.Range("C2").Value = 9.69444257126632
$ws.Range("D2").Value = 5.999117726759839
$ws.Range("E2").Value = 11.79695219597465
$ws.Range("G2").Value = 3.658881143502812
$ws.Range("L2").Value = 9.905606205952498
$ws.Range("M2").Value = 15.717485889696
$ws.Range("N2").Value = 18.92308446137799
$ws.Range("O2").Value = 26.4926244725041

$ws.Range("B3").Value = 16.6922219281363
$ws.Range("C3").Value = 9.442203274051638
$ws.Range("D3").Value = 5.883223368905773
$ws.Range("E3").Value = 11.82643943626278
$ws.Range("G3").Value = 3.661493853868047
$ws.Range("L3").Value = 9.914286092087561
$ws.Range("M3").Value = 15.64014634852755
$ws.Range("N3").Value = 18.99127745930816
$ws.Range("O3").Value = 26.51080969034299

$ws.Range("B4").Value = 16.43795873322783
$ws.Range("C4").Value = 9.282282746470626
$ws.Range("D4").Value = 5.812736658206038
$ws.Range("E4").Value = 11.84553356296851
$ws.Range("G4").Value = 3.663183388450464
$ws.Range("L4").Value = 9.920992661896937
$ws.Range("M4").Value = 15.59526896086543
$ws.Range("N4").Value = 19.03506673561865
$ws.Range("O4").Value = 26.52889550351179

$ws.Range("B5").Value = 16.33409742804043
$ws.Range("C5").Value = 9.215901330781652
$ws.Range("D5").Value = 5.784225573019732
$ws.Range("E5").Value = 11.85356386815422
$ws.Range("G5").Value = 3.663893413655773
$ws.Range("L5").Value = 9.924072273665455
$ws.Range("M5").Value = 15.57765121416274
$ws.Range("N5").Value = 19.05339522506157
$ws.Range("O5").Value = 26.53800210475507

$ws.Range("B6").Value = 16.31684117954007
$ws.Range("C6").Value = 9.204807233628234
$ws.Range("D6").Value = 5.779505461442646
$ws.Range("E6").Value = 11.85491236981903
$ws.Range("G6").Value = 3.664012615036808
$ws.Range("L6").Value = 9.924604585472171
$ws.Range("M6").Value = 15.57476666830076
$ws.Range("N6").Value = 19.05646793592733
$ws.Range("O6").Value = 26.53961900514472

$ws.Range("B7").Value = 16.43655880908109
$ws.Range("C7").Value = 9.281392335783885
$ws.Range("D7").Value = 5.81235122823278
$ws.Range("E7").Value = 11.84564085223492
$ws.Range("G7").Value = 3.663192876844765
$ws.Range("L7").Value = 9.921032790709452
$ws.Range("M7").Value = 15.5950286305188
$ws.Range("N7").Value = 19.03531195824483
$ws.Range("O7").Value = 26.52901129285792

$ws.Range("B8").Value = 16.96225067646162
$ws.Range("C8").Value = 9.608549757234936
$ws.Range("D8").Value = 5.959043765720309
$ws.Range("E8").Value = 11.8069145489879
$ws.Range("G8").Value = 3.659764341171195
$ws.Range("L8").Value = 9.908313444358027
$ws.Range("M8").Value = 15.69028628809487
$ws.Range("N8").Value = 18.9462001781922
$ws.Range("O8").Value = 26.49745658196353

$ws.Range("B9").Value = 17.97368190471378
$ws.Range("C9").Value = 10.20783322272616
$ws.Range("D9").Value = 6.250153324930147
$ws.Range("E9").Value = 11.73878905624889
$ws.Range("G9").Value = 3.65371469599857
$ws.Range("L9").Value = 9.894280390614203
$ws.Range("M9").Value = 15.89713806916071
$ws.Range("N9").Value = 18.7866014650318
$ws.Range("O9").Value = 26.49060850561526

$ws.Range("B10").Value = 18.69657519228935
$ws.Range("C10").Value = 10.61964789636252
$ws.Range("D10").Value = 6.463660387667329
$ws.Range("E10").Value = 11.6934613078288
$ws.Range("G10").Value = 3.649676153421995
$ws.Range("L10").Value = 9.890595819694225
$ws.Range("M10").Value = 16.06041136148213
$ws.Range("N10").Value = 18.67847565771861
$ws.Range("O10").Value = 26.51923788231234

$ws.Range("B11").Value = 19.01945679543922
$ws.Range("C11").Value = 10.80031694677004
$ws.Range("D11").Value = 6.560213135321558
$ws.Range("E11").Value = 11.67385758482728
$ws.Range("G11").Value = 3.647926129516237
$ws.Range("L11").Value = 9.890351555386955
$ws.Range("M11").Value = 16.13693083112425
$ws.Range("N11").Value = 18.63124701625287
$ws.Range("O11").Value = 26.53957470504876

$ws.Range("B12").Value = 19.14074540529089
$ws.Range("C12").Value = 10.86773673022886
$ws.Range("D12").Value = 6.596653820271857
$ws.Range("E12").Value = 11.66657962671245
$ws.Range("G12").Value = 3.647275895929582
$ws.Range("L12").Value = 9.890464245331859
$ws.Range("M12").Value = 16.16621094899051
$ws.Range("N12").Value = 18.61364271494418
$ws.Range("O12").Value = 26.54832580082267

$ws.Range("B13").Value = 19.11466914064866
$ws.Range("C13").Value = 10.85326150492842
$ws.Range("D13").Value = 6.588811686722846
$ws.Range("E13").Value = 11.66814060096688
$ws.Range("G13").Value = 3.64741538214973
$ws.Range("L13").Value = 9.890430861932995
$ws.Range("M13").Value = 16.15989174627833
$ws.Range("N13").Value = 18.61742167888774
$ws.Range("O13").Value = 26.54639443503661

$ws.Range("B14").Value = 19.02945557953474
$ws.Range("C14").Value = 10.80588376594915
$ws.Range("D14").Value = 6.563213790684043
$ws.Range("E14").Value = 11.67325590909436
$ws.Range("G14").Value = 3.647872385021252
$ws.Range("L14").Value = 9.890356718548377
$ws.Range("M14").Value = 16.13933374522717
$ws.Range("N14").Value = 18.62979309319649
$ws.Range("O14").Value = 26.54027363296424

$ws.Range("B15").Value = 18.97712878291487
$ws.Range("C15").Value = 10.77673283506112
$ws.Range("D15").Value = 6.547517351296964
$ws.Range("E15").Value = 11.67640812072634
$ws.Range("G15").Value = 3.648153933489912
$ws.Range("L15").Value = 9.89033800259365
$ws.Range("M15").Value = 16.12678037289605
$ws.Range("N15").Value = 18.63740738412385
$ws.Range("O15").Value = 26.53666113666115

$ws.Range("B16").Value = 18.675343445068
$ws.Range("C16").Value = 10.60770334529958
$ws.Range("D16").Value = 6.457335500816262
$ws.Range("E16").Value = 11.69476285507479
$ws.Range("G16").Value = 3.649792269161063
$ws.Range("L16").Value = 9.890640536533821
$ws.Range("M16").Value = 16.05545420372146
$ws.Range("N16").Value = 18.68160144328762
$ws.Range("O16").Value = 26.51805586678874

$ws.Range("B17").Value = 18.48859023434958
$ws.Range("C17").Value = 10.50227410998334
$ws.Range("D17").Value = 6.40183651092467
$ws.Range("E17").Value = 11.70628273137053
$ws.Range("G17").Value = 3.650819603406873
$ws.Range("L17").Value = 9.891192383367859
$ws.Range("M17").Value = 16.0122596657861
$ws.Range("N17").Value = 18.70921361023623
$ws.Range("O17").Value = 26.50851438545007

$ws.Range("B18").Value = 18.38061984087542
$ws.Range("C18").Value = 10.44100839122568
$ws.Range("D18").Value = 6.369863289272266
$ws.Range("E18").Value = 11.71300433292168
$ws.Range("G18").Value = 3.651418703873078
$ws.Range("L18").Value = 9.891644612504445
$ws.Range("M18").Value = 15.9876278542568
$ws.Range("N18").Value = 18.72527981121949
$ws.Range("O18").Value = 26.50371489422331

$ws.Range("B19").Value = 18.34397148209846
$ws.Range("C19").Value = 10.42015860606357
$ws.Range("D19").Value = 6.359030064735614
$ws.Range("E19").Value = 11.71529660109427
$ws.Range("G19").Value = 3.651622960247801
$ws.Range("L19").Value = 9.891820905705474
$ws.Range("M19").Value = 15.97932502344959
$ws.Range("N19").Value = 18.73075126817873
$ws.Range("O19").Value = 26.50220815499373

$ws.Range("B20").Value = 18.50852877271676
$ws.Range("C20").Value = 10.51356225921567
$ws.Range("D20").Value = 6.407750112695211
$ws.Range("E20").Value = 11.70504652316679
$ws.Range("G20").Value = 3.650709393221716
$ws.Range("L20").Value = 9.891119689739627
$ws.Range("M20").Value = 16.01683593707024
$ws.Range("N20").Value = 18.70625517237284
$ws.Range("O20").Value = 26.50945884202068

$ws.Range("B21").Value = 19.05451237853622
$ws.Range("C21").Value = 10.81982705441409
$ws.Range("D21").Value = 6.570736121127386
$ws.Range("E21").Value = 11.67174947350459
$ws.Range("G21").Value = 3.647737814500283
$ws.Range("L21").Value = 9.890372933588395
$ws.Range("M21").Value = 16.14536403622429
$ws.Range("N21").Value = 18.62615171691391
$ws.Range("O21").Value = 26.54204298307006

$ws.Range("B22").Value = 19.40558429179667
$ws.Range("C22").Value = 11.01417070818443
$ws.Range("D22").Value = 6.676530483238598
$ws.Range("E22").Value = 11.65083602333768
$ws.Range("G22").Value = 3.645868326349254
$ws.Range("L22").Value = 9.89108049039984
$ws.Range("M22").Value = 16.23112651700484
$ws.Range("N22").Value = 18.5754318044992
$ws.Range("O22").Value = 26.56945718721253

$ws.Range("B23").Value = 19.21877541334906
$ws.Range("C23").Value = 10.91098945822415
$ws.Range("D23").Value = 6.620144984774676
$ws.Range("E23").Value = 11.66192051016034
$ws.Range("G23").Value = 3.646859485397859
$ws.Range("L23").Value = 9.890593711400689
$ws.Range("M23").Value = 16.18519866744385
$ws.Range("N23").Value = 18.60235309266371
$ws.Range("O23").Value = 26.55426667859977

$ws.Range("B24").Value = 18.499516429852
$ws.Range("C24").Value = 10.50846091828521
$ws.Range("D24").Value = 6.405076776272438
$ws.Range("E24").Value = 11.70560510550632
$ws.Range("G24").Value = 3.650759192850381
$ws.Range("L24").Value = 9.89115213407265
$ws.Range("M24").Value = 16.01476637631308
$ws.Range("N24").Value = 18.70759208516215
$ws.Range("O24").Value = 26.50902971606635

$ws.Range("B25").Value = 17.70303987267696
$ws.Range("C25").Value = 10.05053216604858
$ws.Range("D25").Value = 6.171291586994531
$ws.Range("E25").Value = 11.75638626044031
$ws.Range("G25").Value = 3.655279635503327
$ws.Range("L25").Value = 9.896911084331537
$ws.Range("M25").Value = 15.83912507696009
$ws.Range("N25").Value = 18.82816611510312
$ws.Range("O25").Value = 26.48655410936829
